$d = $word.ActiveDocument

# 1. Title text: "Results for attempt 21" -> "Results for attempt 4"
$d.Content.Find.Execute("Results for attempt 21", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Results for attempt 4", 2)

# 2. Locate the "Quantitative" Heading 1 paragraph and insert two new paragraphs before it:
#    - "Exercise Name: Wall Slides" (Heading 2)
#    - the long exercise-description paragraph with manual line breaks (Normal)
$quantPara = $null
foreach ($pp in $d.Paragraphs) {
    if ($pp.Style.NameLocal -eq "Heading 1" -and $pp.Range.Text.TrimEnd([char]13) -eq "Quantitative") {
        $quantPara = $pp
        break
    }
}
$quantPara.Range.InsertParagraphBefore()
$quantPara.Range.InsertParagraphBefore()

# Re-scan to find the two freshly inserted blank paragraphs (immediately before "Quantitative")
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $pp = $d.Paragraphs($i)
    if ($pp.Style.NameLocal -eq "Heading 1" -and $pp.Range.Text.TrimEnd([char]13) -eq "Quantitative") {
        $titleIdx = $i - 2
        $bodyIdx = $i - 1
        break
    }
}

$titlePara = $d.Paragraphs($titleIdx)
$titlePara.Range.Text = "Exercise Name: Wall Slides"
$titlePara.Style = "Heading 2"

$bodyPara = $d.Paragraphs($bodyIdx)
$nl = [char]11
$bodyText = "Exercise Name: Sit next to a table with your elbow supported just below shoulder height on a rolled up towel." + $nl + `
    "Now make a gentle fist, keep your elbow bent and then rotate your forearm to point upwards." + $nl + `
    "Return to the start position and relax. Make sure you sit up tall whilst you do this." + $nl + `
    "When it is easy for you to do this you can add a light weight " + [char]8211 + " start with half kilo or a small 500ml water bottle." + $nl + `
    "As the exercise gets easier you can increase the weight:" + $nl + `
    [char]8226 + " First to 1 kilo" + $nl + `
    [char]8226 + " Then to 1 and a half kilos" + $nl + `
    [char]8226 + " Then to 2 kilos" + $nl + `
    "Other Tips:" + $nl + `
    "Be patient! It may be 6 -12 weeks before you see a big change in your pain so you need to stick with it." + $nl
$bodyPara.Range.Text = $bodyText
$bodyPara.Style = "Normal"

# 3. Change heading styles for "Quantitative" and "Qualitative" from Heading 1 -> Heading 2
foreach ($pp in $d.Paragraphs) {
    if ($pp.Style.NameLocal -eq "Heading 1" -and ($pp.Range.Text.TrimEnd([char]13) -eq "Quantitative" -or $pp.Range.Text.TrimEnd([char]13) -eq "Qualitative")) {
        $pp.Style = "Heading 2"
    }
}

# 4. Rebuild Table 1 (originally 3 columns: Exercise Name | Repetitions | Duration) into the
#    new 2-column Quantitative table: Repetitions | Duration, data row: 0 | 1.729892
$t1 = $d.Tables(1)
$t1.Columns(1).Delete()
$t1.Columns(1).Width = 216
$t1.Columns(2).Width = 216
$t1.Cell(2,1).Range.Text = "0"
$t1.Cell(2,2).Range.Text = "1.729892"

# 5. Rebuild Table 2 (originally 2 columns: Exercise Name | Accuracy) into the
#    new 1-column Qualitative table: Accuracy, data row: 0.0
$t2 = $d.Tables(2)
$t2.Columns(1).Delete()
$t2.Columns(1).Width = 432
$t2.Cell(2,1).Range.Text = "0.0"

Write-Output "done"
